{"js": "const replacements = [\n  [\"479\u00d76=\", \"896\u00d79=\"],\n  [\"279\u00d73=\", \"635\u00d72=\"],\n  [\"551\u00d78=\", \"578\u00d77=\"],\n  [\"951\u00d75=\", \"684\u00d77=\"],\n  [\"963\u00d76=\", \"688\u00d78=\"],\n  [\"707\u00d72=\", \"993\u00d72=\"],\n  [\"391\u00d79=\", \"412\u00d72=\"],\n  [\"363\u00d79=\", \"758\u00d74=\"],\n  [\"719\u00d76=\", \"169\u00d73=\"],\n  [\"352\u00d77=\", \"539\u00d76=\"],\n  [\"570\u00d72=\", \"923\u00d72=\"],\n  [\"681\u00d78=\", \"107\u00d77=\"],\n  [\"461\u00d78=\", \"686\u00d73=\"],\n  [\"320\u00d77=\", \"400\u00d78=\"],\n  [\"306\u00d77=\", \"538\u00d79=\"],\n  [\"596\u00d78=\", \"636\u00d77=\"],\n  [\"293\u00d76=\", \"551\u00d75=\"],\n  [\"999\u00d75=\", \"413\u00d75=\"],\n  [\"448\u00d74=\", \"908\u00d73=\"],\n  [\"655\u00d72=\", \"508\u00d75=\"],\n  [\"909\u00d78=\", \"855\u00d79=\"],\n  [\"194\u00d76=\", \"237\u00d77=\"],\n  [\"165\u00d78=\", \"385\u00d72=\"],\n  [\"469\u00d73=\", \"898\u00d78=\"],\n  [\"435\u00d77=\", \"173\u00d73=\"],\n];\n\n// The worksheet is a table of \"AAA\u00d7B=\" multiplication prompts; each cell's\n// text is unique in the document, so an exact-text search/replace safely\n// targets the single matching run per pair.\nfor (const [findText, replaceText] of replacements) {\n  const results = context.document.body.search(findText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items,text\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replaceText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# The document is a table of \"AAA\u00d7B=\" multiplication prompts. Each prompt's\n# text is unique across the document, so a plain Find/Replace on the exact\n# \"AAA\u00d7B=\" string safely targets exactly one run per pair, regardless of\n# which table cell it lives in.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Find = \"479\u00d76=\"; Replace = \"896\u00d79=\"},\n    @{Find = \"279\u00d73=\"; Replace = \"635\u00d72=\"},\n    @{Find = \"551\u00d78=\"; Replace = \"578\u00d77=\"},\n    @{Find = \"951\u00d75=\"; Replace = \"684\u00d77=\"},\n    @{Find = \"963\u00d76=\"; Replace = \"688\u00d78=\"},\n    @{Find = \"707\u00d72=\"; Replace = \"993\u00d72=\"},\n    @{Find = \"391\u00d79=\"; Replace = \"412\u00d72=\"},\n    @{Find = \"363\u00d79=\"; Replace = \"758\u00d74=\"},\n    @{Find = \"719\u00d76=\"; Replace = \"169\u00d73=\"},\n    @{Find = \"352\u00d77=\"; Replace = \"539\u00d76=\"},\n    @{Find = \"570\u00d72=\"; Replace = \"923\u00d72=\"},\n    @{Find = \"681\u00d78=\"; Replace = \"107\u00d77=\"},\n    @{Find = \"461\u00d78=\"; Replace = \"686\u00d73=\"},\n    @{Find = \"320\u00d77=\"; Replace = \"400\u00d78=\"},\n    @{Find = \"306\u00d77=\"; Replace = \"538\u00d79=\"},\n    @{Find = \"596\u00d78=\"; Replace = \"636\u00d77=\"},\n    @{Find = \"293\u00d76=\"; Replace = \"551\u00d75=\"},\n    @{Find = \"999\u00d75=\"; Replace = \"413\u00d75=\"},\n    @{Find = \"448\u00d74=\"; Replace = \"908\u00d73=\"},\n    @{Find = \"655\u00d72=\"; Replace = \"508\u00d75=\"},\n    @{Find = \"909\u00d78=\"; Replace = \"855\u00d79=\"},\n    @{Find = \"194\u00d76=\"; Replace = \"237\u00d77=\"},\n    @{Find = \"165\u00d78=\"; Replace = \"385\u00d72=\"},\n    @{Find = \"469\u00d73=\"; Replace = \"898\u00d78=\"},\n    @{Find = \"435\u00d77=\"; Replace = \"173\u00d73=\"},\n)\n\nforeach ($r in $replacements) {\n    $range = $d.Content\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #         MatchAllWordForms, Forward, Wrap(=wdFindContinue), Format,\n    #         ReplaceWith, Replace(=wdReplaceOne))\n    $range.Find.Execute($r.Find, $false, $false, $false, $false, $false, $true, 1, $false, $r.Replace, 2) | Out-Null\n}\n"}
